$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto price/volume figures (including two coin-row swaps:
# rows 20/21 -> ShibaInu/Litecoin, rows 44/45/46 -> Cronos/TrustWalletToken/FTXToken)
# to match the latest scrape.
#
# Column D ("Price") cells are stored as text in the workbook even when the
# value looks numeric (e.g. "60.32"). Assigning a bare numeric-looking string
# through COM Automation would auto-convert the cell to a Number, so a
# leading apostrophe forces text entry; Style is reset to "Normal" afterwards
# so no stray text-format style is left on the cell.

$c = $ws.Range("D2")
$c.Value = "'41.727.20"
$c.Style = "Normal"
$ws.Range("E2").Value = '  +0.12%  '
$c = $ws.Range("D3")
$c.Value = "'2.205.81"
$c.Style = "Normal"
$ws.Range("E3").Value = '  -2.05%  '
$ws.Range("E4").Value = '  +0.00%  '
$c = $ws.Range("D5")
$c.Value = "'229.60"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -2.01%  '
$c = $ws.Range("D6")
$c.Value = "'0.618"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -3.73%  '
$c = $ws.Range("D7")
$c.Value = "'60.28"
$c.Style = "Normal"
$ws.Range("E7").Value = '  -5.59%  '
$ws.Range("E8").Value = '  -0.04%  '
$c = $ws.Range("D9")
$c.Value = "'0.401"
$c.Style = "Normal"
$ws.Range("E9").Value = '  -2.33%  '
$c = $ws.Range("D10")
$c.Value = "'57.27"
$c.Style = "Normal"
$ws.Range("E10").Value = '  -3.93%  '
$c = $ws.Range("D11")
$c.Value = "'0.0885"
$c.Style = "Normal"
$ws.Range("E11").Value = '  -1.24%  '
$ws.Range("E12").Value = '  -1.53%  '
$c = $ws.Range("D13")
$c.Value = "'2.534.23"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.95%  '
$c = $ws.Range("D14")
$c.Value = "'15.39"
$c.Style = "Normal"
$ws.Range("E14").Value = '  -4.28%  '
$c = $ws.Range("D15")
$c.Value = "'22.13"
$c.Style = "Normal"
$ws.Range("E15").Value = '  -4.22%  '
$c = $ws.Range("D16")
$c.Value = "'5.59"
$c.Style = "Normal"
$ws.Range("E16").Value = '  -2.29%  '
$c = $ws.Range("D17")
$c.Value = "'0.793"
$c.Style = "Normal"
$ws.Range("E17").Value = '  -3.77%  '
$c = $ws.Range("D18")
$c.Value = "'2.209.64"
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.96%  '
$c = $ws.Range("D19")
$c.Value = "'41.665.89"
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.46%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range("D20")
$c.Value = "'0.0₃0900"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -3.37%  '
$ws.Range("B21").Value = 'Litecoin'
$ws.Range("C21").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range("D21")
$c.Value = "'71.97"
$c.Style = "Normal"
$ws.Range("E21").Value = '  -3.42%  '
$c = $ws.Range("D22")
$c.Value = "'6.03"
$c.Style = "Normal"
$ws.Range("E22").Value = '  -2.20%  '
$c = $ws.Range("D23")
$c.Value = "'241.96"
$c.Style = "Normal"
$ws.Range("E23").Value = '  -3.90%  '
$c = $ws.Range("D24")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -0.16%  '
$c = $ws.Range("D25")
$c.Value = "'2.36"
$c.Style = "Normal"
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("E26").Value = '  -3.17%  '
$c = $ws.Range("D27")
$c.Value = "'9.65"
$c.Style = "Normal"
$ws.Range("E27").Value = '  -2.09%  '
$c = $ws.Range("D28")
$c.Value = "'168.75"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("E29").Value = '  -5.58%  '
$c = $ws.Range("D30")
$c.Value = "'1.47"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +1.73%  '
$c = $ws.Range("D31")
$c.Value = "'19.74"
$c.Style = "Normal"
$ws.Range("E31").Value = '  -3.56%  '
$c = $ws.Range("D32")
$c.Value = "'2.62"
$c.Style = "Normal"
$ws.Range("E32").Value = '  -7.65%  '
$ws.Range("E33").Value = '  -2.89%  '
$c = $ws.Range("D34")
$c.Value = "'5.04"
$c.Style = "Normal"
$ws.Range("E34").Value = '  -2.46%  '
$c = $ws.Range("D35")
$c.Value = "'4.62"
$c.Style = "Normal"
$ws.Range("E35").Value = '  -3.62%  '
$c = $ws.Range("D36")
$c.Value = "'0.0646"
$c.Style = "Normal"
$ws.Range("E36").Value = '  +1.30%  '
$c = $ws.Range("D37")
$c.Value = "'2.36"
$c.Style = "Normal"
$ws.Range("E37").Value = '  -4.00%  '
$c = $ws.Range("D38")
$c.Value = "'6.31"
$c.Style = "Normal"
$ws.Range("E38").Value = '  -7.82%  '
$c = $ws.Range("D39")
$c.Value = "'3.54"
$c.Style = "Normal"
$ws.Range("E39").Value = '  -8.66%  '
$c = $ws.Range("D40")
$c.Value = "'0.000238"
$c.Style = "Normal"
$ws.Range("E40").Value = '  -9.76%  '
$c = $ws.Range("D41")
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Range("E41").Value = '  +0.14%  '
$c = $ws.Range("D42")
$c.Value = "'0.0238"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -1.40%  '
$c = $ws.Range("D43")
$c.Value = "'8.58"
$c.Style = "Normal"
$ws.Range("B44").Value = 'Cronos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c = $ws.Range("D44")
$c.Value = "'0.0955"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -3.33%  '
$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$c = $ws.Range("D45")
$c.Value = "'1.20"
$c.Style = "Normal"
$ws.Range("E45").Value = '  -3.12%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Range("D46")
$c.Value = "'4.40"
$c.Style = "Normal"
$ws.Range("E46").Value = '  -13.81%  '
$c = $ws.Range("D47")
$c.Value = "'97.15"
$c.Style = "Normal"
$ws.Range("E47").Value = '  -4.90%  '
$c = $ws.Range("D48")
$c.Value = "'1.466.13"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -2.65%  '
$c = $ws.Range("D49")
$c.Value = "'16.22"
$c.Style = "Normal"
$ws.Range("E49").Value = '  -8.36%  '
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("E51").Value = '  -4.84%  '
